# Adds the "NXsource" / "NXentry" / "NXsample"(name) static fields and the
# supporting "kind" / "static_value" columns required by the NXtomo
# application definition, per the commit message:
#   "added source and other fields necessary in the nxtomo application class"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Stamp formatting on the brand-new cells first (copy format from a
#    neighbouring cell of the same row/kind), so the later .Value writes do
#    not leave them with the engine's "no style" default.
# ---------------------------------------------------------------------------

# Header row (style s=3, like existing I1/J1 header cells) -> K1:L1
$ws.Range("I1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)

# Plain data-row style (s=1, like existing G2) -> new blank/value cells
$ws.Range("G2").Copy()
$ws.Range("K2:L2").PasteSpecial(-4122)

$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("K3:L3").PasteSpecial(-4122)

$ws.Range("K4:L4").PasteSpecial(-4122)
$ws.Range("K5:L5").PasteSpecial(-4122)
$ws.Range("K6:L6").PasteSpecial(-4122)

$ws.Range("I7:L7").PasteSpecial(-4122)

$ws.Range("I8:J8").PasteSpecial(-4122)
$ws.Range("L8").PasteSpecial(-4122)

$ws.Range("A9:B9").PasteSpecial(-4122)
$ws.Range("I9:J9").PasteSpecial(-4122)
$ws.Range("L9").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Header row: rename/extend the trailing columns.
#    J1 was "custom_field" (the old "no" flag column); it becomes
#    "static_value". New columns K ("custom_field") and L ("kind") appended.
# ---------------------------------------------------------------------------
$ws.Range("J1").Value = "static_value"
$ws.Range("K1").Value = "custom_field"
$ws.Range("L1").Value = "kind"

# ---------------------------------------------------------------------------
# 3) Existing data rows (2-6): drop the old "no" flag in J, mark as "group"
#    in the new "kind" column L. Clear() (not ClearContents()) so the cell
#    is removed outright, matching the source's deleted <c> elements.
# ---------------------------------------------------------------------------
$ws.Range("J2:J6").Clear()
$ws.Range("L2").Value = "group"
$ws.Range("L3").Value = "group"
$ws.Range("L4").Value = "group"
$ws.Range("L5").Value = "group"
$ws.Range("L6").Value = "group"

# Row 3 (NXmonitor/control) gains an odin_topic/f142/uint32 module + data_name
$ws.Range("C3").Value = "odin_topic"
$ws.Range("E3").Value = "f142"
$ws.Range("F3").Value = "uint32"
$ws.Range("I3").Value = "data"

# ---------------------------------------------------------------------------
# 4) New row 7: NXsource / source, with a static "probe" = "x-ray" field.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "NXsource"
$ws.Range("B7").Value = "source"
$ws.Range("F7").Value = "string"
$ws.Range("I7").Value = "probe"
$ws.Range("J7").Value = "x-ray"
$ws.Range("L7").Value = "static_data"

# ---------------------------------------------------------------------------
# 5) New row 8: NXentry / entry, with a static "definition" = "NXtomo" field.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "NXentry"
$ws.Range("B8").Value = "entry"
$ws.Range("F8").Value = "string"
$ws.Range("I8").Value = "definition"
$ws.Range("J8").Value = "NXtomo"
$ws.Range("L8").Value = "static_data"

# ---------------------------------------------------------------------------
# 6) New row 9: NXsample / sample, with a static "name" = "lego" field.
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "NXsample"
$ws.Range("B9").Value = "sample"
$ws.Range("F9").Value = "string"
$ws.Range("I9").Value = "name"
$ws.Range("J9").Value = "lego"
$ws.Range("L9").Value = "static_data"

# ---------------------------------------------------------------------------
# 7) Cosmetic bits that accompanied the edit: selection moved to E3, and the
#    sheet picked up an explicit (A4 portrait) page setup.
# ---------------------------------------------------------------------------
[void]$ws.Range("E3").Select()

$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait
